# Add DNA_SampleID column (P) to the "Metadata" sheet (the big per-timepoint
# data log) for every row whose DNA_Sample flag is TRUE, then restore the
# view state (active sheet / selection) the author left the workbook in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Header
$ws.Range("P1").Value = "DNA_SampleID"

# One DNA sample id per bottle/timepoint row that has DNA_Sample = TRUE
$ws.Range("P2").Value  = "144_A0_S6"
$ws.Range("P6").Value  = "144_A4_S7"
$ws.Range("P10").Value = "144_A8_S8"
$ws.Range("P12").Value = "144_B0_S9"
$ws.Range("P16").Value = "144_B4_S10"
$ws.Range("P20").Value = "144_B8_S11"
$ws.Range("P22").Value = "144_C0_S12"
$ws.Range("P26").Value = "144_C4_S13"
$ws.Range("P30").Value = "144_C8_S14"
$ws.Range("P32").Value = "144_D0_S15"
$ws.Range("P36").Value = "144_D4_S16"
$ws.Range("P40").Value = "144_D8_S17"
$ws.Range("P42").Value = "144_E0_S18"
$ws.Range("P46").Value = "144_E4_S19"
$ws.Range("P50").Value = "144_E8_S20"
$ws.Range("P52").Value = "144_F0_S21"
$ws.Range("P56").Value = "144_F4_S22"
$ws.Range("P60").Value = "144_F8_S23"
$ws.Range("P62").Value = "144_G0_S24"
$ws.Range("P66").Value = "144_G4_S25"
$ws.Range("P70").Value = "144_G8_S26"
$ws.Range("P72").Value = "144_H0_S27"
$ws.Range("P76").Value = "144_H4_S28"
$ws.Range("P80").Value = "144_H8_S29"

# The author left the workbook on the "Metadata" tab, scrolled down near the
# bottom of the table with G69 selected (instead of "Data" w/ D19 selected).
$ws.Activate()
[void]$excel.Goto($ws.Range("A54"), $true)
[void]$ws.Range("G69").Select()
